$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The codeforiati:group-name (column C) and codeforiati:group-code (column D)
# values are swapped for every row of the table, including the header row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cValue = $cCell.Value2
    $dValue = $dCell.Value2
    $cCell.Value2 = $dValue
    $dCell.Value2 = $cValue
}
